# The workbook tracks weekly "Berenjena" (eggplant) price records for
# "Femacal de La Calera". This edit adds one more weekly record, inserted
# as a new row 62 (pushing the existing rows 62-165 down to 63-166).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new blank row at position 62; rows 62..165 move to 63..166.
$ws.Rows.Item(62).Insert()

# Populate the newly inserted row 62 with the new weekly record.
$ws.Cells.Item(62, 1).Value = 3
$ws.Cells.Item(62, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(62, 3).Value = "Coquimbo"
$ws.Cells.Item(62, 4).Value = 44495
$ws.Cells.Item(62, 5).Value = 5
$ws.Cells.Item(62, 6).Value = 100112001
$ws.Cells.Item(62, 7).Value = "Berenjena"
$ws.Cells.Item(62, 8).Value = "Sin especificar"
$ws.Cells.Item(62, 9).Value = "Primera"
$ws.Cells.Item(62, 10).Value = 80
$ws.Cells.Item(62, 11).Value = 8500
$ws.Cells.Item(62, 12).Value = 9000
$ws.Cells.Item(62, 13).Value = 8750
$ws.Cells.Item(62, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(62, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(62, 16).Value = 146
$ws.Cells.Item(62, 17).Value = 60
$ws.Cells.Item(62, 18).Value = "Hortaliza"

# Match the date-formatted style used by the rest of column D (sibling cell).
$ws.Cells.Item(62, 4).NumberFormat = $ws.Cells.Item(63, 4).NumberFormat
